# Applies the diff: appends candidate rows 113-132 to the 'Candidates' sheet
# and leaves the dimension to auto-expand to A1:AG132.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 113
$ws.Range("A113").Value = "'2026-01-09 16:51:24"
$ws.Range("B113").Value = "'JAGADEESH M"

# Row 114
$ws.Range("A114").Value = "'2026-01-09 16:57:10"
$ws.Range("B114").Value = "'Test Candidate API"

# Row 115
$ws.Range("A115").Value = "'2026-01-09 17:02:02"
$ws.Range("B115").Value = "'Test Candidate"

# Row 116
$ws.Range("A116").Value = "'2026-01-09 17:02:23"
$ws.Range("B116").Value = "'Test Candidate with File"
$ws.Range("F116").Value = "'1767958343_dummy_resume.pdf"

# Row 117
$ws.Range("A117").Value = "'2026-01-09 17:03:01"
$ws.Range("B117").Value = "'Test Candidate with File"
$ws.Range("F117").Value = "'1767958381_dummy_resume.pdf"

# Row 118
$ws.Range("A118").Value = "'2026-01-09 17:03:40"
$ws.Range("B118").Value = "'Test Candidate with File"
$ws.Range("F118").Value = "'1767958420_dummy_resume.pdf"

# Row 119
$ws.Range("A119").Value = "'2026-01-09 17:09:32"

# Row 120
$ws.Range("A120").Value = "'2026-01-09 17:10:32"
$ws.Range("B120").Value = "'API Test User 2"

# Row 121
$ws.Range("A121").Value = "'2026-01-09 17:11:19"
$ws.Range("B121").Value = "'API Test User Final"

# Row 122
$ws.Range("A122").Value = "'2026-01-09 17:12:17"
$ws.Range("B122").Value = "'API Test User Unique"

# Row 123
$ws.Range("A123").Value = "'2026-01-09 17:14:04"
$ws.Range("B123").Value = "'API Test User Final"

# Row 124
$ws.Range("A124").Value = "'2026-01-09 17:14:34"
$ws.Range("B124").Value = "'API Test Manual"
$ws.Range("C124").Value = "'manual_1767959073784@example.com"
$ws.Range("D124").Value = "'2222222222"
$ws.Range("G124").Value = "'DevOps Engineer"

# Row 125
$ws.Range("A125").Value = "'2026-01-09 17:18:50"
$ws.Range("B125").Value = "'Fixed API Test"

# Row 126
$ws.Range("A126").Value = "'2026-01-09 17:19:55"
$ws.Range("B126").Value = "'Fixed API Test"

# Row 127
$ws.Range("A127").Value = "'2026-01-09 17:23:15"
$ws.Range("B127").Value = "'Final Test"

# Row 128
$ws.Range("A128").Value = "'2026-01-09 17:25:26"
$ws.Range("B128").Value = "'Observed Test"

# Row 129
$ws.Range("A129").Value = "'2026-01-09 17:27:36"
$ws.Range("B129").Value = "'UniqueName12345"

# Row 130
$ws.Range("A130").Value = "'2026-01-09 17:34:16"
$ws.Range("B130").Value = "'Final Test Candidate"

# Row 131
$ws.Range("A131").Value = "'2026-01-09 17:35:11"
$ws.Range("B131").Value = "'Second Test Candidate"
$ws.Range("F131").Value = "'1767960311_dummy_resume.pdf"

# Row 132
$ws.Range("A132").Value = "'2026-01-09 17:35:41"
$ws.Range("B132").Value = "'Final Test Candidate"
$ws.Range("F132").Value = "'1767960341_resume.pdf"
